$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1 — paragraph "Agreagar lat and log to the bat2021_v2 data base."
#
# The trailing run gets split so "data base" is flagged as a grammar error
# (mirrors Word's "data base" -> "database" grammar suggestion):
#   " and log to the bat2021_v2 data base."
#     -> " and log to the bat2021_v2 " + [gramStart]"data base"[gramEnd] + "."
# ---------------------------------------------------------------------------
$target1 = "Agreagar lat and log to the bat2021_v2 data base."
$p1 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13) -eq $target1) {
        $p1 = $para
        break
    }
}
if ($p1 -eq $null) { $p1 = $d.Paragraphs(2) }

$xml1 = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body>' +
                  '<w:p>' +
                    '<w:proofErr w:type="spellStart"/>' +
                    '<w:r><w:t>Agreagar</w:t></w:r>' +
                    '<w:proofErr w:type="spellEnd"/>' +
                    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
                    '<w:proofErr w:type="spellStart"/>' +
                    '<w:r><w:t>lat</w:t></w:r>' +
                    '<w:proofErr w:type="spellEnd"/>' +
                    '<w:r><w:t xml:space="preserve"> and log to the bat2021_v2 </w:t></w:r>' +
                    '<w:proofErr w:type="gramStart"/>' +
                    '<w:r><w:t>data base</w:t></w:r>' +
                    '<w:proofErr w:type="gramEnd"/>' +
                    '<w:r><w:t>.</w:t></w:r>' +
                  '</w:p>' +
                  '<w:p/>' +
                '</w:body>' +
              '</w:document>' +
            '</pkg:xmlData>' +
          '</pkg:part>' +
        '</pkg:package>'

$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2 — paragraph "I tried copying the lat and log ..."
#
# * appends the "Update: I tried again ..." / "However ..." / "I tried the
#   package raster ..." sentences to the existing paragraph;
# * moves "For the population analysis ..." (with its leading line break)
#   into its own paragraph, right after that;
# * adds a brand-new paragraph about calculating riparian vegetation with
#   NDVI;
# * adds one new trailing blank paragraph.
# ---------------------------------------------------------------------------
$target2 = "I tried copying the lat and log but there seems to be a problem as points are not plotting properly. There might be problems with the coordinates. " + [char]11 + "For the population analysis I want to include moon phase, moon illumination, elevation, riparian veg."
$p2 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13) -eq $target2) {
        $p2 = $para
        break
    }
}
if ($p2 -eq $null) { $p2 = $d.Paragraphs(5) }

$xml2 = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body>' +
                  '<w:p>' +
                    '<w:r><w:t xml:space="preserve">I tried copying the </w:t></w:r>' +
                    '<w:proofErr w:type="spellStart"/>' +
                    '<w:r><w:t>lat</w:t></w:r>' +
                    '<w:proofErr w:type="spellEnd"/>' +
                    '<w:r><w:t xml:space="preserve"> and log but there seems to be a problem as points are not plotting properly. There might be problems with the coordinates. </w:t></w:r>' +
                    '<w:r><w:t xml:space="preserve">Update: I tried again and the </w:t></w:r>' +
                    '<w:proofErr w:type="spellStart"/>' +
                    '<w:r><w:t>elevatr</w:t></w:r>' +
                    '<w:proofErr w:type="spellEnd"/>' +
                    '<w:r><w:t xml:space="preserve"> package loaded and worked. </w:t></w:r>' +
                    '<w:proofErr w:type="gramStart"/>' +
                    '<w:r><w:t>However</w:t></w:r>' +
                    '<w:proofErr w:type="gramEnd"/>' +
                    '<w:r><w:t xml:space="preserve"> I didn' + [char]8217 + 't displayed any elevation. </w:t></w:r>' +
                    '<w:r><w:br/><w:t xml:space="preserve">I tried the package raster and plotted the DEM for Idaho and was able to obtain elevation for each site. </w:t></w:r>' +
                  '</w:p>' +
                  '<w:p>' +
                    '<w:r><w:br/><w:t>For the population analysis I want to include moon phase, moon illumination, elevation, riparian veg.</w:t></w:r>' +
                  '</w:p>' +
                  '<w:p>' +
                    '<w:r><w:t xml:space="preserve">The riparian vegetation I will calculate with the </w:t></w:r>' +
                    '<w:r><w:t xml:space="preserve">NDVI. Unfortunately, is not working as it is. I am </w:t></w:r>' +
                    '<w:proofErr w:type="gramStart"/>' +
                    '<w:r><w:t>using  a</w:t></w:r>' +
                    '<w:proofErr w:type="gramEnd"/>' +
                    '<w:r><w:t xml:space="preserve"> code that is in the cleanup script. ChatGPT help me generate it but still not working well. </w:t></w:r>' +
                  '</w:p>' +
                  '<w:p/>' +
                  '<w:p/>' +
                '</w:body>' +
              '</w:document>' +
            '</pkg:xmlData>' +
          '</pkg:part>' +
        '</pkg:package>'

$p2.Range.InsertXML($xml2)

Write-Output "done"
